{"js": "// The document contains a single table of two-digit \u00f7 one-digit division\n// problems. Only every 4th row (0, 4, 8, 12, 16) actually holds text; the\n// rows in between are blank spacer rows. We overwrite the 25 filled cells\n// (5 rows x 5 columns), in row-major / reading order, with their new\n// answer strings - matching the order the diff lists the replacements in.\nconst newValues = [\n  [\"71\u00f73=23, 2\", \"41\u00f78=5, 1\", \"57\u00f75=11, 2\", \"28\u00f74=7, 0\", \"19\u00f79=2, 1\"],\n  [\"86\u00f79=9, 5\", \"95\u00f77=13, 4\", \"28\u00f76=4, 4\", \"42\u00f78=5, 2\", \"80\u00f72=40, 0\"],\n  [\"54\u00f79=6, 0\", \"22\u00f72=11, 0\", \"88\u00f77=12, 4\", \"76\u00f77=10, 6\", \"15\u00f74=3, 3\"],\n  [\"96\u00f79=10, 6\", \"59\u00f75=11, 4\", \"81\u00f78=10, 1\", \"77\u00f72=38, 1\", \"60\u00f77=8, 4\"],\n  [\"83\u00f74=20, 3\", \"96\u00f79=10, 6\", \"55\u00f79=6, 1\", \"58\u00f75=11, 3\", \"67\u00f79=7, 4\"],\n];\n\nconst table = context.document.body.tables.getFirst();\ntable.load(\"values\");\nawait context.sync();\n\n// Detect which rows actually contain text (rather than assuming a fixed\n// stride) so this keeps working even if the spacer-row pattern changes.\nconst filledRowIndexes = [];\ntable.values.forEach((rowValues, idx) => {\n  if (rowValues.some((cellText) => cellText.trim().length > 0)) {\n    filledRowIndexes.push(idx);\n  }\n});\n\nfor (let i = 0; i < filledRowIndexes.length && i < newValues.length; i++) {\n  const rowIndex = filledRowIndexes[i];\n  const rowValues = newValues[i];\n  for (let col = 0; col < rowValues.length; col++) {\n    const cell = table.getCell(rowIndex, col);\n    const paragraph = cell.body.paragraphs.getFirst();\n    const range = paragraph.getRange();\n    // Replace just the run's text content in place so the existing run\n    // formatting (font, size) and paragraph formatting (alignment) are\n    // preserved - only the <w:t> value itself changes, as in the diff.\n    range.insertText(rowValues[col], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document contains a single table of two-digit \u00f7 one-digit division\n# problems. Only every 4th row (1, 5, 9, 13, 17 in Word's 1-based Rows\n# collection) actually holds text; the rows in between are blank spacer\n# rows. We overwrite the 25 filled cells (5 rows x 5 columns), in\n# row-major / reading order, with their new answer strings - matching the\n# order the diff lists the replacements in.\n\n$newValues = @(\n    @(\"71\u00f73=23, 2\", \"41\u00f78=5, 1\", \"57\u00f75=11, 2\", \"28\u00f74=7, 0\", \"19\u00f79=2, 1\"),\n    @(\"86\u00f79=9, 5\", \"95\u00f77=13, 4\", \"28\u00f76=4, 4\", \"42\u00f78=5, 2\", \"80\u00f72=40, 0\"),\n    @(\"54\u00f79=6, 0\", \"22\u00f72=11, 0\", \"88\u00f77=12, 4\", \"76\u00f77=10, 6\", \"15\u00f74=3, 3\"),\n    @(\"96\u00f79=10, 6\", \"59\u00f75=11, 4\", \"81\u00f78=10, 1\", \"77\u00f72=38, 1\", \"60\u00f77=8, 4\"),\n    @(\"83\u00f74=20, 3\", \"96\u00f79=10, 6\", \"55\u00f79=6, 1\", \"58\u00f75=11, 3\", \"67\u00f79=7, 4\")\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Detect which rows actually contain text (rather than assuming a fixed\n# stride) so this keeps working even if the spacer-row pattern changes.\n$filledRowIndexes = @()\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    $hasText = $false\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        $txt = $t.Cell($r, $c).Range.Text\n        $txt = $txt.TrimEnd([char]13, [char]7)\n        if ($txt.Trim().Length -gt 0) {\n            $hasText = $true\n        }\n    }\n    if ($hasText) {\n        $filledRowIndexes += $r\n    }\n}\n\nfor ($i = 0; $i -lt $filledRowIndexes.Length -and $i -lt $newValues.Length; $i++) {\n    $rowIndex = $filledRowIndexes[$i]\n    $rowValues = $newValues[$i]\n    for ($col = 1; $col -le $rowValues.Length; $col++) {\n        $cell = $t.Cell($rowIndex, $col)\n        # Assigning to Range.Text (rather than replacing the whole cell)\n        # keeps the run's existing formatting (font, size) and the\n        # paragraph's alignment intact - only the <w:t> value changes,\n        # just like in the diff.\n        $cell.Range.Text = $rowValues[$col - 1]\n    }\n}\n"}
